# Realestate Update resale numbers 2025-01-07 08:53
# Append a new data row (row 12) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Text / string columns (A-D). B (time) and C (weekday) are plain words
# that Excel never reinterprets, so a normal assignment keeps them as text.
# A (date) and D (week number) look like a date / number to Excel, so they
# need a leading apostrophe to force text entry; the style is then reset
# back to Normal so the cell doesn't keep the "quote prefix" formatting.
$ws.Cells.Item($row, 1).Value = "'2025-01-07"
$ws.Cells.Item($row, 2).Value = "08:53:02"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "'01"

$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 4).Style = "Normal"

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 127378
$ws.Cells.Item($row, 6).Value = 143554
$ws.Cells.Item($row, 7).Value = 168905
$ws.Cells.Item($row, 8).Value = 158469
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 141987
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192449
$ws.Cells.Item($row, 14).Value = 114975
$ws.Cells.Item($row, 15).Value = 45585
$ws.Cells.Item($row, 16).Value = 28325
$ws.Cells.Item($row, 17).Value = 63921
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47080
$ws.Cells.Item($row, 20).Value = -1
